$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# Header row (row 1)
$ws.Range("B1").Value = "bank"
$ws.Range("C1").Value = "deposit_type"
$ws.Range("D1").Value = "currency"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

# Data rows (2-21)
# Row 2
$ws.Range("B2").Value = "立法院郵局（第25支局）"
$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "潘维剛"
$ws.Range("F2").Value = 1054022
$ws.Range("G2").Value = "deposit"
$ws.Range("H2").Value = "normal"
$ws.Range("I2").Value = "2012-03-26"
$ws.Range("J2").Value = "潘維剛"
$ws.Range("K2").Value = 678
$ws.Range("L2").Value = "tmp71a01"
$ws.Range("M2").Value = 45

# Row 3
$ws.Range("B3").Value = "台北長安郵局（第46支局）"
$ws.Range("C3").Value = "其他存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "潘维剛"
$ws.Range("F3").Value = 2251
$ws.Range("G3").Value = "deposit"
$ws.Range("H3").Value = "normal"
$ws.Range("I3").Value = "2012-03-26"
$ws.Range("J3").Value = "潘維剛"
$ws.Range("K3").Value = 678
$ws.Range("L3").Value = "tmp71a01"
$ws.Range("M3").Value = 46

# Row 4
$ws.Range("B4").Value = "永豐商業銀行板新分行"
$ws.Range("C4").Value = "活期儲蓄存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "田正超"
$ws.Range("F4").Value = 494765
$ws.Range("G4").Value = "deposit"
$ws.Range("H4").Value = "normal"
$ws.Range("I4").Value = "2012-03-26"
$ws.Range("J4").Value = "潘維剛"
$ws.Range("K4").Value = 678
$ws.Range("L4").Value = "tmp71a01"
$ws.Range("M4").Value = 47

# Row 5
$ws.Range("B5").Value = "永豐商業銀行西松分行"
$ws.Range("C5").Value = "活期儲蓄存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "田正超"
$ws.Range("F5").Value = 1372036
$ws.Range("G5").Value = "deposit"
$ws.Range("H5").Value = "normal"
$ws.Range("I5").Value = "2012-03-26"
$ws.Range("J5").Value = "潘維剛"
$ws.Range("K5").Value = 678
$ws.Range("L5").Value = "tmp71a01"
$ws.Range("M5").Value = 48

# Row 6
$ws.Range("B6").Value = "日盛國際商業銀行松山分行"
$ws.Range("C6").Value = "活期儲蓄存款"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("E6").Value = "田正超"
$ws.Range("F6").Value = 2374298
$ws.Range("G6").Value = "deposit"
$ws.Range("H6").Value = "normal"
$ws.Range("I6").Value = "2012-03-26"
$ws.Range("J6").Value = "潘維剛"
$ws.Range("K6").Value = 678
$ws.Range("L6").Value = "tmp71a01"
$ws.Range("M6").Value = 49

# Row 7
$ws.Range("B7").Value = "彰化商業銀行總行"
$ws.Range("C7").Value = "活期儲蓄存款"
$ws.Range("D7").Value = "新臺幣"
$ws.Range("E7").Value = "田正超"
$ws.Range("F7").Value = 2314717
$ws.Range("G7").Value = "deposit"
$ws.Range("H7").Value = "normal"
$ws.Range("I7").Value = "2012-03-26"
$ws.Range("J7").Value = "潘維剛"
$ws.Range("K7").Value = 678
$ws.Range("L7").Value = "tmp71a01"
$ws.Range("M7").Value = 50

# Row 8
$ws.Range("B8").Value = "中國信託商業銀行城中分行"
$ws.Range("C8").Value = "活期儲蓄存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "田正超"
$ws.Range("F8").Value = 679609
$ws.Range("G8").Value = "deposit"
$ws.Range("H8").Value = "normal"
$ws.Range("I8").Value = "2012-03-26"
$ws.Range("J8").Value = "潘維剛"
$ws.Range("K8").Value = 678
$ws.Range("L8").Value = "tmp71a01"
$ws.Range("M8").Value = 51

# Row 9
$ws.Range("B9").Value = "三信商業銀行台中分行"
$ws.Range("C9").Value = "活期儲蓄存款"
$ws.Range("D9").Value = "新臺幣"
$ws.Range("E9").Value = "田正超"
$ws.Range("F9").Value = 5282
$ws.Range("G9").Value = "deposit"
$ws.Range("H9").Value = "normal"
$ws.Range("I9").Value = "2012-03-26"
$ws.Range("J9").Value = "潘維剛"
$ws.Range("K9").Value = 678
$ws.Range("L9").Value = "tmp71a01"
$ws.Range("M9").Value = 52

# Row 10
$ws.Range("B10").Value = "台新國際商業銀行敦北分行"
$ws.Range("C10").Value = "活期儲蓄存款"
$ws.Range("D10").Value = "新臺幣"
$ws.Range("E10").Value = "田正超"
$ws.Range("F10").Value = 1266328
$ws.Range("G10").Value = "deposit"
$ws.Range("H10").Value = "normal"
$ws.Range("I10").Value = "2012-03-26"
$ws.Range("J10").Value = "潘維剛"
$ws.Range("K10").Value = 678
$ws.Range("L10").Value = "tmp71a01"
$ws.Range("M10").Value = 53

# Row 11
$ws.Range("B11").Value = "花旗（台灣）銀行營業部"
$ws.Range("C11").Value = "活期儲蓄存款"
$ws.Range("D11").Value = "新臺幣"
$ws.Range("E11").Value = "田正超"
$ws.Range("F11").Value = 1147.16
$ws.Range("G11").Value = "deposit"
$ws.Range("H11").Value = "normal"
$ws.Range("I11").Value = "2012-03-26"
$ws.Range("J11").Value = "潘維剛"
$ws.Range("K11").Value = 678
$ws.Range("L11").Value = "tmp71a01"
$ws.Range("M11").Value = 54

# Row 12
$ws.Range("B12").Value = "台新國際商業銀行營業部"
$ws.Range("C12").Value = "活期儲蓄存款"
$ws.Range("D12").Value = "新臺幣"
$ws.Range("E12").Value = "潘維剛"
$ws.Range("F12").Value = 40706
$ws.Range("G12").Value = "deposit"
$ws.Range("H12").Value = "normal"
$ws.Range("I12").Value = "2012-03-26"
$ws.Range("J12").Value = "潘維剛"
$ws.Range("K12").Value = 678
$ws.Range("L12").Value = "tmp71a01"
$ws.Range("M12").Value = 55

# Row 13
$ws.Range("B13").Value = "台新國際商業銀行營業部"
$ws.Range("C13").Value = "綜合存款"
$ws.Range("D13").Value = "美金"
$ws.Range("E13").Value = "潘維剛"
$ws.Range("F13").Value = 8676.13
$ws.Range("G13").Value = "deposit"
$ws.Range("H13").Value = "normal"
$ws.Range("I13").Value = "2012-03-26"
$ws.Range("J13").Value = "潘維剛"
$ws.Range("K13").Value = 678
$ws.Range("L13").Value = "tmp71a01"
$ws.Range("M13").Value = 56

# Row 14
$ws.Range("B14").Value = "花旗（台灣）銀行營業部"
$ws.Range("C14").Value = "活期存款"
$ws.Range("D14").Value = "新臺幣"
$ws.Range("E14").Value = "潘維剛"
$ws.Range("F14").Value = 101048.04
$ws.Range("G14").Value = "deposit"
$ws.Range("H14").Value = "normal"
$ws.Range("I14").Value = "2012-03-26"
$ws.Range("J14").Value = "潘維剛"
$ws.Range("K14").Value = 678
$ws.Range("L14").Value = "tmp71a01"
$ws.Range("M14").Value = 57

# Row 15
$ws.Range("B15").Value = "花旗（台灣）銀行營業部"
$ws.Range("C15").Value = "活期存款"
$ws.Range("D15").Value = "美金"
$ws.Range("E15").Value = "潘維剛"
$ws.Range("F15").Value = 366.96
$ws.Range("G15").Value = "deposit"
$ws.Range("H15").Value = "normal"
$ws.Range("I15").Value = "2012-03-26"
$ws.Range("J15").Value = "潘維剛"
$ws.Range("K15").Value = 678
$ws.Range("L15").Value = "tmp71a01"
$ws.Range("M15").Value = 58

# Row 16
$ws.Range("B16").Value = "台北富邦商業銀行敦南分行"
$ws.Range("C16").Value = "活期存款"
$ws.Range("D16").Value = "新臺幣"
$ws.Range("E16").Value = "潘維剛"
$ws.Range("F16").Value = 1151980.57
$ws.Range("G16").Value = "deposit"
$ws.Range("H16").Value = "normal"
$ws.Range("I16").Value = "2012-03-26"
$ws.Range("J16").Value = "潘維剛"
$ws.Range("K16").Value = 678
$ws.Range("L16").Value = "tmp71a01"
$ws.Range("M16").Value = 59

# Row 17
$ws.Range("B17").Value = "玉山商業銀行民生分行"
$ws.Range("C17").Value = "活期存款"
$ws.Range("D17").Value = "新臺幣"
$ws.Range("E17").Value = "潘維剛"
$ws.Range("F17").Value = 5517
$ws.Range("G17").Value = "deposit"
$ws.Range("H17").Value = "normal"
$ws.Range("I17").Value = "2012-03-26"
$ws.Range("J17").Value = "潘維剛"
$ws.Range("K17").Value = 678
$ws.Range("L17").Value = "tmp71a01"
$ws.Range("M17").Value = 60

# Row 18
$ws.Range("B18").Value = "永豐商業銀行松江分行"
$ws.Range("C18").Value = "活期儲蓄存款"
$ws.Range("D18").Value = "新毫幣"
$ws.Range("E18").Value = "潘維剛"
$ws.Range("F18").Value = 1178566
$ws.Range("G18").Value = "deposit"
$ws.Range("H18").Value = "normal"
$ws.Range("I18").Value = "2012-03-26"
$ws.Range("J18").Value = "潘維剛"
$ws.Range("K18").Value = 678
$ws.Range("L18").Value = "tmp71a01"
$ws.Range("M18").Value = 61

# Row 19
$ws.Range("B19").Value = "永豐商業銀行松江分行"
$ws.Range("C19").Value = "綜合存款"
$ws.Range("D19").Value = "美金"
$ws.Range("E19").Value = "潘維剛"
$ws.Range("F19").Value = 7239414.63
$ws.Range("G19").Value = "deposit"
$ws.Range("H19").Value = "normal"
$ws.Range("I19").Value = "2012-03-26"
$ws.Range("J19").Value = "潘維剛"
$ws.Range("K19").Value = 678
$ws.Range("L19").Value = "tmp71a01"
$ws.Range("M19").Value = 62

# Row 20
$ws.Range("B20").Value = "國泰世華商業銀行南京東路分行"
$ws.Range("C20").Value = "活期儲蓄存款"
$ws.Range("D20").Value = "新臺幣"
$ws.Range("E20").Value = "潘維剛"
$ws.Range("F20").Value = 1552
$ws.Range("G20").Value = "deposit"
$ws.Range("H20").Value = "normal"
$ws.Range("I20").Value = "2012-03-26"
$ws.Range("J20").Value = "潘維剛"
$ws.Range("K20").Value = 678
$ws.Range("L20").Value = "tmp71a01"
$ws.Range("M20").Value = 63

# Row 21
$ws.Range("B21").Value = "台新國際商業銀行營業部"
$ws.Range("C21").Value = "活期儲蓄存款"
$ws.Range("D21").Value = "新臺幣"
$ws.Range("E21").Value = "潘維剛"
$ws.Range("F21").Value = 40706
$ws.Range("G21").Value = "deposit"
$ws.Range("H21").Value = "normal"
$ws.Range("I21").Value = "2012-03-26"
$ws.Range("J21").Value = "潘維剛"
$ws.Range("K21").Value = 678
$ws.Range("L21").Value = "tmp71a01"
$ws.Range("M21").Value = 64
